# Quarterly indexing esoteric bug-fix operation
#
# The sheet holds, per base-quarter row (A2:A16 = Q10..Q24), a run of QoQ
# naive-forecaster errors in columns B..K. Those errors were off-by-one in
# the quarterly index: the error that belongs in the first slot (column B,
# "1 quarter ahead") was missing, and every existing error was sitting one
# column too far to the left. The fix: insert the correct "1-quarter-ahead"
# error into column B for each row and shift everything that was already
# there one column to the right (B->C, C->D, ... ), letting any value that
# would fall past column K (the 10th slot) drop off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The newly-computed "1 quarter ahead" error to insert into column B of each row.
$newFirstQuarterError = @{
    2  = -1.025188112727922
    3  = 0.08364543516793629
    4  = -0.1538585523806955
    5  = 0.7495351060200912
    6  = 0.03849281619118239
    7  = -0.2590580299438133
    8  = 0.01855976243503714
    9  = 0.1467044301255134
    10 = -0.1819613811903656
    11 = 0.4718454808444464
    12 = -0.08594117411414147
    13 = -0.07695400962807622
    14 = -0.5068991247689255
    15 = 0.6215838649243215
    16 = -0.2766911554241067
}

$colB = 2   # column B, first data slot
$colK = 11  # column K, last data slot (10th slot)

for ($row = 2; $row -le 16; $row++) {

    # Find the last populated data column in this row (starting the scan at B).
    $lastCol = $colB - 1
    for ($col = $colB; $col -le $colK; $col++) {
        if ($ws.Cells.Item($row, $col).Value2 -ne $null) {
            $lastCol = $col
        }
    }

    # Shift existing values one column to the right, working from the
    # rightmost value back towards B so nothing gets clobbered before it is
    # read. Anything shifting past column K simply falls off the end.
    for ($col = $lastCol; $col -ge $colB; $col--) {
        $value = $ws.Cells.Item($row, $col).Value2
        $destCol = $col + 1
        if ($destCol -le $colK) {
            $ws.Cells.Item($row, $destCol).Value = $value
        }
    }

    # Insert the newly-computed first-quarter-ahead error into column B.
    $ws.Cells.Item($row, $colB).Value = $newFirstQuarterError[$row]
}
